$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

# Update the confidential disclosure date string (shared string used in A16)
$ws.Range("A16").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-04 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-13
$ws.Range("D2").Value = 0.03097120733000712
$ws.Range("E2").Value = -0.005754908598510355

$ws.Range("D3").Value = 0.02390218175336856
$ws.Range("E3").Value = -0.004864489228630919

$ws.Range("D4").Value = 0.05135613448571712
$ws.Range("E4").Value = -0.004913430042115086

$ws.Range("D5").Value = 0.1377079921208408
$ws.Range("E5").Value = -0.006143896523847947

$ws.Range("D6").Value = 0.03018004003664627
$ws.Range("E6").Value = 0

$ws.Range("D7").Value = 0.1210945391904059
$ws.Range("E7").Value = -0.01050160612799589

$ws.Range("D8").Value = 0.1015707545174934
$ws.Range("E8").Value = 0.002932013927066324

$ws.Range("D9").Value = 0.02835161425340419
$ws.Range("E9").Value = 0.01398601398601396

$ws.Range("D10").Value = 0.1232223032840157
$ws.Range("E10").Value = 0.00630975143403445

$ws.Range("D11").Value = 0.2481562892321967
$ws.Range("E11").Value = -0.02020384233787309

$ws.Range("D12").Value = 0.1034869437959043
$ws.Range("E12").Value = -0.01185695161598777

$ws.Range("E13").Value = -0.007433509139124883
